$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("FAPs", "Tac1", "Tacr2", "ECs", 3, 1, 14.882091, 44.646273, 0.996401763178, 0.996401763178, 2, 0.6666666666666666, 0.2385553333333333, 0.7156659999999999, 0.9144827885830529, 0.914482788583053, 3.550202179202, 31.951819612818, 0.911192262940088, 0.9111922629400881)
for ($i = 0; $i -lt $row2.Count; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$row3 = @("FAPs", "Tac1", "Tacr2", "MuSCs", 3, 1, 14.882091, 44.646273, 0.996401763178, 0.996401763178, 1, 0.3333333333333333, 0.02230833333333333, 0.066925, 0.08551721141694704, 0.08551721141694704, 0.331994646725, 2.987951820525, 0.08520950023791182, 0.08520950023791182)
for ($i = 0; $i -lt $row3.Count; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

$row4 = @("MuSCs", "Tac1", "Tacr2", "ECs", 1, 0.3333333333333333, 0.037615, 0.112845, 0.002518439937098924, 0.002518439937098924, 2, 0.6666666666666666, 0.2385553333333333, 0.7156659999999999, 0.9144827885830529, 0.914482788583053, 0.008973258863333333, 0.08075932976999999, 0.002303069976557153, 0.002303069976557153)
for ($i = 0; $i -lt $row4.Count; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

$row5 = @("MuSCs", "Tac1", "Tacr2", "MuSCs", 1, 0.3333333333333333, 0.037615, 0.112845, 0.002518439937098924, 0.002518439937098924, 1, 0.3333333333333333, 0.02230833333333333, 0.066925, 0.08551721141694704, 0.08551721141694704, 0.0008391279583333333, 0.007552151625, 0.0002153699605417715, 0.0002153699605417715)
for ($i = 0; $i -lt $row5.Count; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

$row6 = @("Resolving-Mac", "Tac1", "Tacr2", "ECs", 1, 0.3333333333333333, 0.01612766666666667, 0.048383, 0.001079796884901035, 0.001079796884901035, 2, 0.6666666666666666, 0.2385553333333333, 0.7156659999999999, 0.9144827885830529, 0.914482788583053, 0.003847340897555556, 0.034626068078, 0.0009874556664075922, 0.000987455666407592)
for ($i = 0; $i -lt $row6.Count; $i++) {
    $ws.Cells.Item(6, $i + 1).Value = $row6[$i]
}

$row7 = @("Resolving-Mac", "Tac1", "Tacr2", "MuSCs", 1, 0.3333333333333333, 0.01612766666666667, 0.048383, 0.001079796884901035, 0.001079796884901035, 1, 0.3333333333333333, 0.02230833333333333, 0.066925, 0.08551721141694704, 0.08551721141694704, 0.0003597813638888889, 0.003238032275, 0.000092341218493442633516964102, 0.00009234121849344260641190979)
for ($i = 0; $i -lt $row7.Count; $i++) {
    $ws.Cells.Item(7, $i + 1).Value = $row7[$i]
}

